# Edit: Added unit req's from paper
# - Inserts a new "Unit definitions" worksheet between "Sets" and "Parameters"
# - Populates it with per-unit scheduling requirement data
# - Restores the prior selections on "Sets" and "Parameters"

$wb = $excel.ActiveWorkbook
$sets = $wb.Worksheets.Item("Sets")
$params = $wb.Worksheets.Item("Parameters")

# Restore Sets' active cell (recorded in workbook as last state before save)
$sets.Activate()
$sets.Range("N18").Select() | Out-Null

# Insert the new sheet right after "Sets" (and thus right before "Parameters")
$ws = $wb.Worksheets.Add($null, $sets)
$ws.Name = "Unit definitions"

# Restore Parameters' active cell before leaving it (it keeps B1 selected,
# it's simply no longer the active/visible tab once we switch away)
$params.Activate()
$params.Range("B1").Select() | Out-Null

$ws.Activate()

# Header: A1 then G1:L1 (R1Min..R3Max), matching original authoring order
$ws.Range("A1").Value = "Unit"
$ws.Range("G1").Value = "R1Min"
$ws.Range("H1").Value = "R1Max"
$ws.Range("I1").Value = "R2Min"
$ws.Range("J1").Value = "R2Max"
$ws.Range("K1").Value = "R3Min"
$ws.Range("L1").Value = "R3Max"

# Data rows 2-19, left to right
$ws.Range("A2").Value = "Geriatrics"
$ws.Range("B2").Value = 2
$ws.Range("C2").Value = 2
$ws.Range("D2").Value = 4
$ws.Range("E2").Value = 4
$ws.Range("F2").Formula = "=FALSE"
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 2
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0

$ws.Range("A3").Value = "R1F"
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 4
$ws.Range("D3").Value = 2
$ws.Range("E3").Value = 4
$ws.Range("F3").Formula = "=TRUE"
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 2
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 1

$ws.Range("A4").Value = "HemeF"
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = 4
$ws.Range("D4").Value = 2
$ws.Range("E4").Value = 4
$ws.Range("F4").Value = $true
$ws.Range("G4").Value = 1
$ws.Range("H4").Value = 1
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0

$ws.Range("A5").Value = "CardF"
$ws.Range("B5").Value = 2
$ws.Range("C5").Value = 4
$ws.Range("D5").Value = 2
$ws.Range("E5").Value = 4
$ws.Range("F5").Value = $true
$ws.Range("G5").Value = 1
$ws.Range("H5").Value = 1
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 1

$ws.Range("A6").Value = "FloatF"
$ws.Range("B6").Value = 2
$ws.Range("C6").Value = 4
$ws.Range("D6").Value = 2
$ws.Range("E6").Value = 4
$ws.Range("F6").Value = $true
$ws.Range("G6").Value = 1
$ws.Range("H6").Value = 1
$ws.Range("I6").Value = 1
$ws.Range("J6").Value = 1
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 0

$ws.Range("A7").Value = "TBC1"
$ws.Range("B7").Value = 2
$ws.Range("C7").Value = 4
$ws.Range("D7").Value = 2
$ws.Range("E7").Value = 4
$ws.Range("F7").Value = $true
$ws.Range("G7").Value = 1
$ws.Range("H7").Value = 1
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 1

$ws.Range("A8").Value = "mat_D"
$ws.Range("B8").Value = 2
$ws.Range("C8").Value = 2
$ws.Range("D8").Value = 2
$ws.Range("E8").Value = "unlimited"
$ws.Range("F8").Value = $true
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 1
$ws.Range("J8").Value = 1
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 1

$ws.Range("A9").Value = "MAT N"
$ws.Range("B9").Value = 2
$ws.Range("C9").Value = 2
$ws.Range("D9").Value = 2
$ws.Range("E9").Value = "unlimited"
$ws.Range("F9").Value = $true
$ws.Range("G9").Value = 0
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 1
$ws.Range("J9").Value = 1
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 1

$ws.Range("A10").Value = "Electives"
$ws.Range("B10").Value = 2
$ws.Range("C10").Value = 4
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = "unlimited"
$ws.Range("F10").Value = $false
$ws.Range("G10").Value = 0
$ws.Range("H10").Value = "unlimited"
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = "unlimited"
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = "unlimited"

$ws.Range("A11").Value = "VAC"
$ws.Range("B11").Value = 2
$ws.Range("C11").Value = 2
$ws.Range("D11").Value = 4
$ws.Range("E11").Value = 4
$ws.Range("F11").Value = $false
$ws.Range("G11").Value = 0
$ws.Range("H11").Value = 2
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 5
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 5

$ws.Range("A12").Value = "RNAT"
$ws.Range("B12").Value = 2
$ws.Range("C12").Value = 4
$ws.Range("D12").Value = 4
$ws.Range("E12").Value = 4
$ws.Range("F12").Value = $true
$ws.Range("G12").Value = 0
$ws.Range("H12").Value = 2
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 1

$ws.Range("A13").Value = "Overnight"
$ws.Range("B13").Value = 2
$ws.Range("C13").Value = 2
$ws.Range("D13").Value = 4
$ws.Range("E13").Value = 4
$ws.Range("F13").Value = $true
$ws.Range("G13").Value = 1
$ws.Range("H13").Value = 2
$ws.Range("I13").Value = 1
$ws.Range("J13").Value = 1
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 0

$ws.Range("A14").Value = "Midcall"
$ws.Range("B14").Value = 2
$ws.Range("C14").Value = 2
$ws.Range("D14").Value = 4
$ws.Range("E14").Value = 4
$ws.Range("F14").Value = $true
$ws.Range("G14").Value = 0
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 1
$ws.Range("J14").Value = 1
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 0

$ws.Range("A15").Value = "Sick call"
$ws.Range("B15").Value = 1
$ws.Range("C15").Value = 2
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = "unlimited"
$ws.Range("F15").Value = $false
$ws.Range("G15").Value = 0
$ws.Range("H15").Value = "unlimited"
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = "unlimited"
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = "unlimited"

$ws.Range("A16").Value = "MICU_D"
$ws.Range("B16").Value = 2
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 8
$ws.Range("F16").Value = $true
$ws.Range("G16").Value = 2
$ws.Range("H16").Value = 2
$ws.Range("I16").Value = 1
$ws.Range("J16").Value = 1
$ws.Range("K16").Value = 1
$ws.Range("L16").Value = 1

$ws.Range("A17").Value = "MICU_N"
$ws.Range("B17").Value = 2
$ws.Range("C17").Value = 2
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = 8
$ws.Range("F17").Value = $true
$ws.Range("G17").Value = 1
$ws.Range("H17").Value = 1
$ws.Range("I17").Value = 1
$ws.Range("J17").Value = 1
$ws.Range("K17").Value = 1
$ws.Range("L17").Value = 1

$ws.Range("A18").Value = "Twig"
$ws.Range("B18").Value = 1
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 10
$ws.Range("E18").Value = 11
$ws.Range("F18").Value = $true
$ws.Range("G18").Value = 1
$ws.Range("H18").Value = 2
$ws.Range("I18").Value = 1
$ws.Range("J18").Value = 2
$ws.Range("K18").Value = 1
$ws.Range("L18").Value = 2

$ws.Range("A19").Value = "OPD"
$ws.Range("B19").Value = 1
$ws.Range("C19").Value = 1
$ws.Range("D19").Value = 10
$ws.Range("E19").Value = 11
$ws.Range("F19").Value = $true
$ws.Range("G19").Value = 1
$ws.Range("H19").Value = 3
$ws.Range("I19").Value = 1
$ws.Range("J19").Value = 2
$ws.Range("K19").Value = 1
$ws.Range("L19").Value = 2

# Header: B1:F1 (Duration/Rotation/Student_Req), added last
$ws.Range("B1").Value = "Duration_Min"
$ws.Range("C1").Value = "Duration_Max"
$ws.Range("D1").Value = "Rotation_Min"
$ws.Range("E1").Value = "Rotation_Max"
$ws.Range("F1").Value = "Student_Req"

# Select C20 on the new sheet to match the final recorded selection
# (this also leaves "Unit definitions" as the active/visible tab)
$ws.Range("C20").Select() | Out-Null
